$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 unchanged text, but shared-string rewrite ("IW Checker" -> "IWChecker")
$ws.Range("A2").Value = "IWChecker"

# Row 3 ("IW Maker" -> "IWMaker")
$ws.Range("A3").Value = "IWMaker"

# Row 4 ("OCM Checker" -> "OCMChecker"), and its "Last Changed By"/"Last Changed On" values change
$ws.Range("A4").Value = "OCMChecker"
$ws.Range("C4").Value = "02/06/2021 15:55:37"

# Row 5 ("OCM Maker" -> "OCMMaker"), "Last Changed On" changes
$ws.Range("A5").Value = "OCMMaker"
$ws.Range("C5").Value = "03/06/2021 15:13:14"

# "Last Changed By" column: all users renamed from EC2AMAZ-J9C4H3C to EC2AMAZ-N8SAHHO\Administrator
$ws.Range("B2").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("B3").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("B4").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("B5").Value = "EC2AMAZ-N8SAHHO\Administrator"

# "Last Changed On" for rows 2 and 3 also updated
$ws.Range("C2").Value = "19/05/2021 06:31:27"
$ws.Range("C3").Value = "19/05/2021 06:30:55"
